# Update cryptos list values (prices and 1h volume %) to latest scrape
# Commit: Updated cryptos list on Fri Jun  7 03:18:02 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.773.02"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").Value = "3.796.89"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'702.49"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D6").Value = "'169.86"
$ws.Range("E6").Value = "  -2.49%  "
$ws.Range("D7").Value = "3.794.78"
$ws.Range("E7").Value = "  -1.87%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.521"
$ws.Range("E9").Value = "  -0.88%  "
$ws.Range("E10").Value = "  -1.80%  "
$ws.Range("D11").Value = "'7.56"
$ws.Range("E11").Value = "  +5.30%  "
$ws.Range("E12").Value = "  -0.46%  "
$ws.Range("E13").Value = "  -4.18%  "
$ws.Range("D14").Value = "'35.68"
$ws.Range("E14").Value = "  -2.22%  "
$ws.Range("D15").Value = "4.436.20"
$ws.Range("E15").Value = "  -1.90%  "
$ws.Range("D16").Value = "3.786.12"
$ws.Range("E16").Value = "  -2.32%  "
$ws.Range("D17").Value = "70.729.08"
$ws.Range("D18").Value = "'0.114"
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'17.34"
$ws.Range("E19").Value = "  -1.95%  "
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").Value = "'7.09"
$ws.Range("E20").Value = "  -2.06%  "
$ws.Range("D21").Value = "'497.10"
$ws.Range("E21").Value = "  -0.41%  "
$ws.Range("D22").Value = "'10.63"
$ws.Range("E22").Value = "  -4.63%  "
$ws.Range("D23").Value = "'0.720"
$ws.Range("E23").Value = "  -0.48%  "
$ws.Range("D24").Value = "'84.09"
$ws.Range("E24").Value = "  -1.04%  "
$ws.Range("E25").Value = "  -4.92%  "
$ws.Range("D26").Value = "3.947.03"
$ws.Range("E26").Value = "  -1.68%  "
$ws.Range("D27").Value = "'12.01"
$ws.Range("E27").Value = "  -2.09%  "
$ws.Range("D28").Value = "'10.27"
$ws.Range("E28").Value = "  -5.13%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("E30").Value = "  -6.94%  "
$ws.Range("D31").Value = "'3.02"
$ws.Range("E31").Value = "  -5.27%  "
$ws.Range("D32").Value = "'2.25"
$ws.Range("E32").Value = "  -1.24%  "
$ws.Range("E33").Value = "  -3.72%  "
$ws.Range("D34").Value = "'28.93"
$ws.Range("E34").Value = "  -2.68%  "
$ws.Range("E35").Value = "  -3.42%  "
$ws.Range("B36").Value = "Binance-PegBSC-USD"
$ws.Range("C36").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  -0.20%  "
$ws.Range("B37").Value = "RenzoRestakedETH"
$ws.Range("C37").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D37").Value = "3.760.99"
$ws.Range("E37").Value = "  -1.61%  "
$ws.Range("D38").Value = "'9.01"
$ws.Range("E38").Value = "  -2.59%  "
$ws.Range("E39").Value = "  -3.92%  "
$ws.Range("D40").Value = "'2.36"
$ws.Range("E40").Value = "  -1.59%  "
$ws.Range("E41").Value = "  -2.94%  "
$ws.Range("E42").Value = "  -1.62%  "
$ws.Range("D44").Value = "'3.22"
$ws.Range("E44").Value = "  -6.34%  "
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").Value = "'166.55"
$ws.Range("E46").Value = "  +1.72%  "
$ws.Range("D47").Value = "'0.000314"
$ws.Range("E47").Value = "  +0.82%  "
$ws.Range("D48").Value = "'48.90"
$ws.Range("E48").Value = "  -1.05%  "
$ws.Range("D49").Value = "'417.01"
$ws.Range("E49").Value = "  -0.29%  "
$ws.Range("E50").Value = "  -1.21%  "
$ws.Range("E51").Value = "  -3.81%  "
